# Automatische test-sync: 2025-06-23 18:16:50
# Append the new "Offerte voor 500 stuks" email to the Logs sheet and the
# matching roll-up row to the Dashboard sheet's category summary, then grow
# the conditional-formatting ranges and the chart series to cover the new row.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs!A7:G7 -------------------------------------------------------------
$logs.Range("A7").Value = "Offerte voor 500 stuks"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Graag ontvang ik een offerte voor 500 stuks van product X."
$logs.Range("D7").Value = "Offerte / Prijsaanvraag"
$logs.Range("E7").Value = "Beste klant,`nDank u wel voor uw interesse in product X. Om een passende offerte voor u te kunnen opstellen, hebben wij wat meer informatie nodig. Zou u zo vriendelijk willen zijn om de gewenste specificaties van product X (zoals kleur, afmetingen, eventuele personalisatie-opties) met ons te delen? Op basis van deze informatie kunnen wij u een nauwkeurige offerte sturen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F7").Value = "2025-06-23 18:16:02"
$logs.Range("G7").Value = "Ja"

# The multi-line value in E7 makes the host auto-grow the row height; put it
# back to the sheet's default (un-setting the custom-height flag) so row 7
# ends up plain, like the other data rows.
$logs.Rows.Item(7).EntireRow.AutoFit()

# --- Dashboard!A6:B6 ----------------------------------------------------------
$dash.Range("A6").Value = "Offerte / Prijsaanvraag"
$dash.Range("B6").Value = 1

# --- Grow the conditional formatting ranges on Logs to include row 7 --------
# Each block of cfRules shares one sqref; modifying one rule's AppliesTo
# range moves the whole block (dxfIds/priorities/formulas stay untouched).
$logs.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D7"))
$logs.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G7"))

# --- Extend the Dashboard bar chart's category/value series to row 6 --------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$6,Dashboard!`$B`$2:`$B`$6,1)"
